# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append the per-play yardage gains logged during the Wild Card
# round to the running season play logs (Rush/Pass, Offense/Defense).
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 7 -3 9 -5 2 11 6 -1 9 2 14 -4 -3 34 6 5 6"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 1 6 0 3 16 1 32 35 6 0 10 28 2 5 7 31 19 17 4 20 16 11"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 17 3 2 2 2 2 3 2 1 5 6 -1 5 7 -2 2 8 7 -2 7 -5 6 5 2 9 3 3 3 5 1"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 8 2 6 9 3 8 5 13 8 17 16 8 18 9 8 5 13 -4 4 3 22 5 17 2 36 16 3 5 6"

# ---------------------------------------------------------------------------
# OFF sheet: updated offensive situational / drive totals after the game.
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 403
$offWs.Range("E2").Value = 24
$offWs.Range("F2").Value = 134
$offWs.Range("J2").Value = 79
$offWs.Range("N2").Value = 54
$offWs.Range("O2").Value = 60
$offWs.Range("P2").Value = 25

$offWs.Range("B3").Value = 23
$offWs.Range("C3").Value = 352
$offWs.Range("E3").Value = 69
$offWs.Range("F3").Value = 192
$offWs.Range("G3").Value = 69
$offWs.Range("H3").Value = 55
$offWs.Range("I3").Value = 107
$offWs.Range("J3").Value = 110
$offWs.Range("L3").Value = 594
$offWs.Range("M3").Value = 355
$offWs.Range("Q3").Value = 1138

# ---------------------------------------------------------------------------
# DEF sheet: updated defensive situational / drive totals after the game.
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 422
$defWs.Range("D2").Value = 33
$defWs.Range("E2").Value = 22
$defWs.Range("F2").Value = 119
$defWs.Range("G2").Value = 134
$defWs.Range("H2").Value = 13
$defWs.Range("J2").Value = 56
$defWs.Range("N2").Value = 46
$defWs.Range("O2").Value = 44

$defWs.Range("C3").Value = 368
$defWs.Range("D3").Value = 10
$defWs.Range("E3").Value = 72
$defWs.Range("F3").Value = 208
$defWs.Range("G3").Value = 82
$defWs.Range("I3").Value = 114
$defWs.Range("J3").Value = 113
$defWs.Range("L3").Value = 598
$defWs.Range("M3").Value = 423
$defWs.Range("Q3").Value = 1140

# ---------------------------------------------------------------------------
# ST sheet: updated special-teams totals plus appended per-kick/punt logs.
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 161
$stWs.Range("D2").Value = 130
$stWs.Range("F2").Value = 183
$stWs.Range("G2").Value = 174
$stWs.Range("H2").Value = 22
$stWs.Range("I2").Value = 13

$stWs.Range("B3").Value = 104

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 61"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 18"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 17 16 16 24"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 27 36 46 36 58 49"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 0 1 0 8 4"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 5 0 -1 31"

# ---------------------------------------------------------------------------
# TURNS sheet: updated turnover totals.
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 19
$turnsWs.Range("D3").Value = 19

# ---------------------------------------------------------------------------
# PEN sheet: updated penalty totals.
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B3").Value = 28
$penWs.Range("B4").Value = 6
